$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fn1"
$ws.Cells.Item(2, 3).Value = "Mag"
$ws.Cells.Item(2, 4).Value = "M1"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 40.7349555
$ws.Cells.Item(2, 8).Value = 81.469911
$ws.Cells.Item(2, 9).Value = 0.05567871843833241
$ws.Cells.Item(2, 10).Value = 0.03826666865920979
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.27644
$ws.Cells.Item(2, 14).Value = 0.8293199999999999
$ws.Cells.Item(2, 15).Value = 0.1940440920813295
$ws.Cells.Item(2, 16).Value = 0.2208748168298663
$ws.Cells.Item(2, 17).Value = 11.26077109842
$ws.Cells.Item(2, 18).Value = 67.56462659051999
$ws.Cells.Item(2, 19).Value = 0.01080412636761819
$ws.Cells.Item(2, 20).Value = 0.008452143430792147

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fn1"
$ws.Cells.Item(3, 3).Value = "Mag"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 40.7349555
$ws.Cells.Item(3, 8).Value = 81.469911
$ws.Cells.Item(3, 9).Value = 0.05567871843833241
$ws.Cells.Item(3, 10).Value = 0.03826666865920979
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.6290166666666667
$ws.Cells.Item(3, 14).Value = 1.88705
$ws.Cells.Item(3, 15).Value = 0.4415315004607062
$ws.Cells.Item(3, 16).Value = 0.502582625643659
$ws.Cells.Item(3, 17).Value = 25.622965925425
$ws.Cells.Item(3, 18).Value = 153.73779555255
$ws.Cells.Item(3, 19).Value = 0.0245839080958061
$ws.Cells.Item(3, 20).Value = 0.01923216280938157

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fn1"
$ws.Cells.Item(4, 3).Value = "Mag"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 40.7349555
$ws.Cells.Item(4, 8).Value = 81.469911
$ws.Cells.Item(4, 9).Value = 0.05567871843833241
$ws.Cells.Item(4, 10).Value = 0.03826666865920979
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.5191680000000001
$ws.Cells.Item(4, 14).Value = 1.038336
$ws.Cells.Item(4, 15).Value = 0.3644244074579644
$ws.Cells.Item(4, 16).Value = 0.2765425575264748
$ws.Cells.Item(4, 17).Value = 21.148285377024
$ws.Cells.Item(4, 18).Value = 84.59314150809601
$ws.Cells.Item(4, 19).Value = 0.02029068397490813
$ws.Cells.Item(4, 20).Value = 0.01058236241903607

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fn1"
$ws.Cells.Item(5, 3).Value = "Mag"
$ws.Cells.Item(5, 4).Value = "M1"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 350.3919066666667
$ws.Cells.Item(5, 8).Value = 1051.17572
$ws.Cells.Item(5, 9).Value = 0.4789344206933965
$ws.Cells.Item(5, 10).Value = 0.4937404802104949
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.27644
$ws.Cells.Item(5, 14).Value = 0.8293199999999999
$ws.Cells.Item(5, 15).Value = 0.1940440920813295
$ws.Cells.Item(5, 16).Value = 0.2208748168298663
$ws.Cells.Item(5, 17).Value = 96.86233867893333
$ws.Cells.Item(5, 18).Value = 871.7610481103999
$ws.Cells.Item(5, 19).Value = 0.09293439482994761
$ws.Cells.Item(5, 20).Value = 0.1090548381279833

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fn1"
$ws.Cells.Item(6, 3).Value = "Mag"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 350.3919066666667
$ws.Cells.Item(6, 8).Value = 1051.17572
$ws.Cells.Item(6, 9).Value = 0.4789344206933965
$ws.Cells.Item(6, 10).Value = 0.4937404802104949
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.6290166666666667
$ws.Cells.Item(6, 14).Value = 1.88705
$ws.Cells.Item(6, 15).Value = 0.4415315004607062
$ws.Cells.Item(6, 16).Value = 0.502582625643659
$ws.Cells.Item(6, 17).Value = 220.4023491584445
$ws.Cells.Item(6, 18).Value = 1983.621142426
$ws.Cells.Item(6, 19).Value = 0.2114646333910344
$ws.Cells.Item(6, 20).Value = 0.2481453869307516

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fn1"
$ws.Cells.Item(7, 3).Value = "Mag"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 350.3919066666667
$ws.Cells.Item(7, 8).Value = 1051.17572
$ws.Cells.Item(7, 9).Value = 0.4789344206933965
$ws.Cells.Item(7, 10).Value = 0.4937404802104949
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.5191680000000001
$ws.Cells.Item(7, 14).Value = 1.038336
$ws.Cells.Item(7, 15).Value = 0.3644244074579644
$ws.Cells.Item(7, 16).Value = 0.2765425575264748
$ws.Cells.Item(7, 17).Value = 181.91226540032
$ws.Cells.Item(7, 18).Value = 1091.47359240192
$ws.Cells.Item(7, 19).Value = 0.1745353924724145
$ws.Cells.Item(7, 20).Value = 0.1365402551517601

# Row 8
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Fn1"
$ws.Cells.Item(8, 3).Value = "Mag"
$ws.Cells.Item(8, 4).Value = "M1"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 243.8287033333334
$ws.Cells.Item(8, 8).Value = 731.4861100000001
$ws.Cells.Item(8, 9).Value = 0.3332781281688242
$ws.Cells.Item(8, 10).Value = 0.3435812836494235
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.27644
$ws.Cells.Item(8, 14).Value = 0.8293199999999999
$ws.Cells.Item(8, 15).Value = 0.1940440920813295
$ws.Cells.Item(8, 16).Value = 0.2208748168298663
$ws.Cells.Item(8, 17).Value = 67.40400674946666
$ws.Cells.Item(8, 18).Value = 606.6360607452
$ws.Cells.Item(8, 19).Value = 0.06467065179108446
$ws.Cells.Item(8, 20).Value = 0.07588845309223674

# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Fn1"
$ws.Cells.Item(9, 3).Value = "Mag"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 243.8287033333334
$ws.Cells.Item(9, 8).Value = 731.4861100000001
$ws.Cells.Item(9, 9).Value = 0.3332781281688242
$ws.Cells.Item(9, 10).Value = 0.3435812836494235
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.6290166666666667
$ws.Cells.Item(9, 14).Value = 1.88705
$ws.Cells.Item(9, 15).Value = 0.4415315004607062
$ws.Cells.Item(9, 16).Value = 0.502582625643659
$ws.Cells.Item(9, 17).Value = 153.3723182083889
$ws.Cells.Item(9, 18).Value = 1380.3508638755
$ws.Cells.Item(9, 19).Value = 0.1471527920011165
$ws.Cells.Item(9, 20).Value = 0.172677983658546

# Row 10
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Fn1"
$ws.Cells.Item(10, 3).Value = "Mag"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 243.8287033333334
$ws.Cells.Item(10, 8).Value = 731.4861100000001
$ws.Cells.Item(10, 9).Value = 0.3332781281688242
$ws.Cells.Item(10, 10).Value = 0.3435812836494235
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.5191680000000001
$ws.Cells.Item(10, 14).Value = 1.038336
$ws.Cells.Item(10, 15).Value = 0.3644244074579644
$ws.Cells.Item(10, 16).Value = 0.2765425575264748
$ws.Cells.Item(10, 17).Value = 126.58806025216
$ws.Cells.Item(10, 18).Value = 759.5283615129601
$ws.Cells.Item(10, 19).Value = 0.1214546843766233
$ws.Cells.Item(10, 20).Value = 0.09501484689864076

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Fn1"
$ws.Cells.Item(11, 3).Value = "Mag"
$ws.Cells.Item(11, 4).Value = "M1"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 48.737294
$ws.Cells.Item(11, 8).Value = 146.211882
$ws.Cells.Item(11, 9).Value = 0.06661674320651284
$ws.Cells.Item(11, 10).Value = 0.06867617226847689
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.27644
$ws.Cells.Item(11, 14).Value = 0.8293199999999999
$ws.Cells.Item(11, 15).Value = 0.1940440920813295
$ws.Cells.Item(11, 16).Value = 0.2208748168298663
$ws.Cells.Item(11, 17).Value = 13.47293755336
$ws.Cells.Item(11, 18).Value = 121.25643798024
$ws.Cells.Item(11, 19).Value = 0.01292658545292286
$ws.Cells.Item(11, 20).Value = 0.01516883697037618

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Fn1"
$ws.Cells.Item(12, 3).Value = "Mag"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 48.737294
$ws.Cells.Item(12, 8).Value = 146.211882
$ws.Cells.Item(12, 9).Value = 0.06661674320651284
$ws.Cells.Item(12, 10).Value = 0.06867617226847689
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.6290166666666667
$ws.Cells.Item(12, 14).Value = 1.88705
$ws.Cells.Item(12, 15).Value = 0.4415315004607062
$ws.Cells.Item(12, 16).Value = 0.502582625643659
$ws.Cells.Item(12, 17).Value = 30.65657021423333
$ws.Cells.Item(12, 18).Value = 275.9091319281
$ws.Cells.Item(12, 19).Value = 0.02941339058377717
$ws.Cells.Item(12, 20).Value = 0.03451545097784735

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Fn1"
$ws.Cells.Item(13, 3).Value = "Mag"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 48.737294
$ws.Cells.Item(13, 8).Value = 146.211882
$ws.Cells.Item(13, 9).Value = 0.06661674320651284
$ws.Cells.Item(13, 10).Value = 0.06867617226847689
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.5191680000000001
$ws.Cells.Item(13, 14).Value = 1.038336
$ws.Cells.Item(13, 15).Value = 0.3644244074579644
$ws.Cells.Item(13, 16).Value = 0.2765425575264748
$ws.Cells.Item(13, 17).Value = 25.302843451392
$ws.Cells.Item(13, 18).Value = 151.817060708352
$ws.Cells.Item(13, 19).Value = 0.02427676716981282
$ws.Cells.Item(13, 20).Value = 0.01899188432025337

# Row 14
$ws.Cells.Item(14, 1).Value = "Neutro"
$ws.Cells.Item(14, 2).Value = "Fn1"
$ws.Cells.Item(14, 3).Value = "Mag"
$ws.Cells.Item(14, 4).Value = "M1"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 22.832077
$ws.Cells.Item(14, 8).Value = 68.49623099999999
$ws.Cells.Item(14, 9).Value = 0.03120810544755168
$ws.Cells.Item(14, 10).Value = 0.03217289111905
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.27644
$ws.Cells.Item(14, 14).Value = 0.8293199999999999
$ws.Cells.Item(14, 15).Value = 0.1940440920813295
$ws.Cells.Item(14, 16).Value = 0.2208748168298663
$ws.Cells.Item(14, 17).Value = 6.311699365879998
$ws.Cells.Item(14, 18).Value = 56.80529429291999
$ws.Cells.Item(14, 19).Value = 0.006055748487148559
$ws.Cells.Item(14, 20).Value = 0.007106181432807401

# Row 15
$ws.Cells.Item(15, 1).Value = "Neutro"
$ws.Cells.Item(15, 2).Value = "Fn1"
$ws.Cells.Item(15, 3).Value = "Mag"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 22.832077
$ws.Cells.Item(15, 8).Value = 68.49623099999999
$ws.Cells.Item(15, 9).Value = 0.03120810544755168
$ws.Cells.Item(15, 10).Value = 0.03217289111905
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.6290166666666667
$ws.Cells.Item(15, 14).Value = 1.88705
$ws.Cells.Item(15, 15).Value = 0.4415315004607062
$ws.Cells.Item(15, 16).Value = 0.502582625643659
$ws.Cells.Item(15, 17).Value = 14.36175696761667
$ws.Cells.Item(15, 18).Value = 129.25581270855
$ws.Cells.Item(15, 19).Value = 0.01377936162479343
$ws.Cells.Item(15, 20).Value = 0.01616953609315971

# Row 16
$ws.Cells.Item(16, 1).Value = "Neutro"
$ws.Cells.Item(16, 2).Value = "Fn1"
$ws.Cells.Item(16, 3).Value = "Mag"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 22.832077
$ws.Cells.Item(16, 8).Value = 68.49623099999999
$ws.Cells.Item(16, 9).Value = 0.03120810544755168
$ws.Cells.Item(16, 10).Value = 0.03217289111905
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.5191680000000001
$ws.Cells.Item(16, 14).Value = 1.038336
$ws.Cells.Item(16, 15).Value = 0.3644244074579644
$ws.Cells.Item(16, 16).Value = 0.2765425575264748
$ws.Cells.Item(16, 17).Value = 11.853683751936
$ws.Cells.Item(16, 18).Value = 71.12210251161601
$ws.Cells.Item(16, 19).Value = 0.01137299533560969
$ws.Cells.Item(16, 20).Value = 0.008897173593082897

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Fn1"
$ws.Cells.Item(17, 3).Value = "Mag"
$ws.Cells.Item(17, 4).Value = "M1"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 25.082339
$ws.Cells.Item(17, 8).Value = 50.164678
$ws.Cells.Item(17, 9).Value = 0.03428388404538221
$ws.Cells.Item(17, 10).Value = 0.02356250409334498
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.27644
$ws.Cells.Item(17, 14).Value = 0.8293199999999999
$ws.Cells.Item(17, 15).Value = 0.1940440920813295
$ws.Cells.Item(17, 16).Value = 0.2208748168298663
$ws.Cells.Item(17, 17).Value = 6.93376179316
$ws.Cells.Item(17, 18).Value = 41.60257075896
$ws.Cells.Item(17, 19).Value = 0.006652585152607768
$ws.Cells.Item(17, 20).Value = 0.005204363775670546

# Row 18
$ws.Cells.Item(18, 1).Value = "sCs"
$ws.Cells.Item(18, 2).Value = "Fn1"
$ws.Cells.Item(18, 3).Value = "Mag"
$ws.Cells.Item(18, 4).Value = "M2"
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 25.082339
$ws.Cells.Item(18, 8).Value = 50.164678
$ws.Cells.Item(18, 9).Value = 0.03428388404538221
$ws.Cells.Item(18, 10).Value = 0.02356250409334498
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.6290166666666667
$ws.Cells.Item(18, 14).Value = 1.88705
$ws.Cells.Item(18, 15).Value = 0.4415315004607062
$ws.Cells.Item(18, 16).Value = 0.502582625643659
$ws.Cells.Item(18, 17).Value = 15.77720926998333
$ws.Cells.Item(18, 18).Value = 94.66325561990001
$ws.Cells.Item(18, 19).Value = 0.01513741476417847
$ws.Cells.Item(18, 20).Value = 0.01184210517397278

# Row 19
$ws.Cells.Item(19, 1).Value = "sCs"
$ws.Cells.Item(19, 2).Value = "Fn1"
$ws.Cells.Item(19, 3).Value = "Mag"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 25.082339
$ws.Cells.Item(19, 8).Value = 50.164678
$ws.Cells.Item(19, 9).Value = 0.03428388404538221
$ws.Cells.Item(19, 10).Value = 0.02356250409334498
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.5191680000000001
$ws.Cells.Item(19, 14).Value = 1.038336
$ws.Cells.Item(19, 15).Value = 0.3644244074579644
$ws.Cells.Item(19, 16).Value = 0.2765425575264748
$ws.Cells.Item(19, 17).Value = 13.021947773952
$ws.Cells.Item(19, 18).Value = 52.08779109580801
$ws.Cells.Item(19, 19).Value = 0.01249388412859597
$ws.Cells.Item(19, 20).Value = 0.006516035143701651
